$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 21 (D, M, N, O, P, Q, S, T change; L stays "Primera")
$ws.Cells.Item(21, 4).Value = 44516
$ws.Cells.Item(21, 13).Value = 100
$ws.Cells.Item(21, 14).Value = 1900
$ws.Cells.Item(21, 15).Value = 2000
$ws.Cells.Item(21, 16).Value = 1950
$ws.Cells.Item(21, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(21, 19).Value = 1950
$ws.Cells.Item(21, 20).Value = 1

# New row 22
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44516
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100107
$ws.Cells.Item(22, 8).Value = "Otros"
$ws.Cells.Item(22, 9).Value = 100107002
$ws.Cells.Item(22, 10).Value = "Chirimoya"
$ws.Cells.Item(22, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 50
$ws.Cells.Item(22, 14).Value = 1700
$ws.Cells.Item(22, 15).Value = 1700
$ws.Cells.Item(22, 16).Value = 1700
$ws.Cells.Item(22, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(22, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 19).Value = 1700
$ws.Cells.Item(22, 20).Value = 1

# New row 23 (this is essentially the original row 21 data moved down)
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44491
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100107
$ws.Cells.Item(23, 8).Value = "Otros"
$ws.Cells.Item(23, 9).Value = 100107002
$ws.Cells.Item(23, 10).Value = "Chirimoya"
$ws.Cells.Item(23, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 150
$ws.Cells.Item(23, 14).Value = 25000
$ws.Cells.Item(23, 15).Value = 26000
$ws.Cells.Item(23, 16).Value = 25467
$ws.Cells.Item(23, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 19).Value = 2547
$ws.Cells.Item(23, 20).Value = 10

# Copy date cell formatting from D21 (style index 2) to D22/D23
$ws.Range("D21").Copy()
$ws.Range("D22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D23").PasteSpecial(-4122)
